$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Repayment strategy row (row 17) changes from "RBI (India)" to the new
# "Overdue/Due Fee/Int,Principal" scenario value.
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Reflect the author's final selection (the edited cell).
$ws.Range("B17").Select()
